$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formats from column O (rows 3-25) down to the new column P so that
# borders / fonts / number formats line up with the existing table.
$ws.Range("O3:O25").Copy() | Out-Null
$ws.Range("P3:P25").PasteSpecial(-4122) | Out-Null

# A few cells in the new column ended up with slightly different
# (but already-existing) styles than their column-O counterparts.
$ws.Range("O8").Copy() | Out-Null
$ws.Range("P5").PasteSpecial(-4122) | Out-Null

$ws.Range("O6").Copy() | Out-Null
$ws.Range("P7").PasteSpecial(-4122) | Out-Null
$ws.Range("P8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# New 2021 column values
$ws.Range("P4").Value = 2021
$ws.Range("P5").Value = 9038
$ws.Range("P7").Value = 8587
$ws.Range("P8").Value = 451

# Rows without 2021 data yet show the same "no data" placeholder ("…")
# used in column O.
$ws.Range("P10:P25").Value = "…"

# Selection marker left behind by the editor (cell just right of the new column).
$ws.Range("Q4").Select() | Out-Null
